$d = $word.ActiveDocument

# Remove the bold "Whatever I want " run (together with the trailing
# space run) that used to sit between "So i can do" and the underlined
# "Whit it." run.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Whatever I want ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = ""
}

# Remove everything from the leading space/superscript/subscript run
# after "Whit it." through the final "!DOCTYPE html" run (numbered and
# bulleted list lines included), leaving the paragraph ending right
# after "Whit it."
$startRng = $d.Content
$startRng.Find.Execute(" But it has to be correct displayed", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("!DOCTYPE html", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end = $endRng.End

$delRng = $d.Range($start, $end)
$delRng.Text = ""
